$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue $ws.Range("D2") "26.152.72"
Set-TextValue $ws.Range("E2") "  -0.57%  "

# Row 3: Ethereum
Set-TextValue $ws.Range("D3") "1.663.73"
Set-TextValue $ws.Range("E3") "  +0.24%  "

# Row 4: TetherUSD
Set-TextValue $ws.Range("E4") "  -0.38%  "

# Row 5: BNB
Set-TextValue $ws.Range("D5") "217.62"
Set-TextValue $ws.Range("E5") "  -0.86%  "

# Row 6: XRP
Set-TextValue $ws.Range("D6") "0.5243"
Set-TextValue $ws.Range("E6") "  +0.25%  "

# Row 7: USDC
Set-TextValue $ws.Range("E7") "  -0.27%  "

# Row 8: Cardano
Set-TextValue $ws.Range("D8") "0.2638"
Set-TextValue $ws.Range("E8") "  -0.93%  "

# Row 9: Dogecoin
Set-TextValue $ws.Range("D9") "0.06296"
Set-TextValue $ws.Range("E9") "  -0.60%  "

# Row 10: Solana
Set-TextValue $ws.Range("D10") "20.69"
Set-TextValue $ws.Range("E10") "  -3.70%  "

# Row 11: TRON
Set-TextValue $ws.Range("D11") "0.07744"
Set-TextValue $ws.Range("E11") "  -0.16%  "

# Row 12: Polkadot
Set-TextValue $ws.Range("D12") "4.464"
Set-TextValue $ws.Range("E12") "  +0.38%  "

# Row 13: WrappedEther
Set-TextValue $ws.Range("D13") "1.662.03"
Set-TextValue $ws.Range("E13") "  -0.27%  "

# Row 14: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "1.890.98"
Set-TextValue $ws.Range("E14") "  +0.04%  "

# Row 15: Polygon
Set-TextValue $ws.Range("D15") "0.5464"
Set-TextValue $ws.Range("E15") "  -0.51%  "

# Row 16: ShibaInu
Set-TextValue $ws.Range("D16") "0.0₅8124"
Set-TextValue $ws.Range("E16") "  -1.51%  "

# Row 17: Litecoin
Set-TextValue $ws.Range("D17") "64.86"
Set-TextValue $ws.Range("E17") "  -0.50%  "

# Row 18: WrappedBTC
Set-TextValue $ws.Range("D18") "26.162.36"
Set-TextValue $ws.Range("E18") "  -0.66%  "

# Row 19: Dai
Set-TextValue $ws.Range("D19") "1.003"
Set-TextValue $ws.Range("E19") "  -0.26%  "

# Row 20: Uniswap
Set-TextValue $ws.Range("D20") "4.590"
Set-TextValue $ws.Range("E20") "  -2.46%  "

# Row 21: BitcoinCash
Set-TextValue $ws.Range("D21") "191.61"
Set-TextValue $ws.Range("E21") "  -0.32%  "

# Row 22: Avalanche
Set-TextValue $ws.Range("E22") "  -1.98%  "

# Row 23: Chainlink
Set-TextValue $ws.Range("D23") "6.000"
Set-TextValue $ws.Range("E23") "  -3.74%  "

# Row 24: BinanceUSD
Set-TextValue $ws.Range("D24") "1.005"
Set-TextValue $ws.Range("E24") "  -0.43%  "

# Row 25: Monero
Set-TextValue $ws.Range("D25") "137.86"
Set-TextValue $ws.Range("E25") "  -0.73%  "

# Row 26: Stellar
Set-TextValue $ws.Range("D26") "0.1244"
Set-TextValue $ws.Range("E26") "  -1.09%  "

# Row 27: Cosmos
Set-TextValue $ws.Range("D27") "7.252"
Set-TextValue $ws.Range("E27") "  -1.19%  "

# Row 28: EthereumClassic
Set-TextValue $ws.Range("D28") "16.17"
Set-TextValue $ws.Range("E28") "  +0.30%  "

# Row 29: Toncoin
Set-TextValue $ws.Range("D29") "1.403"

# Row 30: Hedera
Set-TextValue $ws.Range("D30") "0.05961"
Set-TextValue $ws.Range("E30") "  -1.92%  "

# Row 31: PancakeSwap
Set-TextValue $ws.Range("D31") "1.279"
Set-TextValue $ws.Range("E31") "  -0.77%  "

# Row 32: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D32") "3.532"
Set-TextValue $ws.Range("E32") "  -0.65%  "

# Row 33: Filecoin
Set-TextValue $ws.Range("D33") "3.266"
Set-TextValue $ws.Range("E33") "  -3.46%  "

# Row 34: LidoDAOToken
Set-TextValue $ws.Range("D34") "1.571"
Set-TextValue $ws.Range("E34") "  -5.85%  "

# Row 35: ARBITRUM
Set-TextValue $ws.Range("D35") "0.9581"
Set-TextValue $ws.Range("E35") "  -3.46%  "

# Row 36: HuobiToken
Set-TextValue $ws.Range("E36") "  -0.39%  "

# Row 37: MXToken
Set-TextValue $ws.Range("D37") "2.772"
Set-TextValue $ws.Range("E37") "  -0.29%  "

# Row 38: ImmutableX
Set-TextValue $ws.Range("D38") "0.5654"
Set-TextValue $ws.Range("E38") "  -5.56%  "

# Row 39: VeChain
Set-TextValue $ws.Range("D39") "0.01595"
Set-TextValue $ws.Range("E39") "  -0.69%  "

# Row 40: FraxShare
Set-TextValue $ws.Range("D40") "5.918"
Set-TextValue $ws.Range("E40") "  -1.20%  "

# Row 41: TrustWalletToken
Set-TextValue $ws.Range("D41") "0.8506"
Set-TextValue $ws.Range("E41") "  -0.67%  "

# Row 42: PaxDollar
Set-TextValue $ws.Range("D42") "1.003"
Set-TextValue $ws.Range("E42") "  -0.19%  "

# Row 43: Quant
Set-TextValue $ws.Range("D43") "101.16"
Set-TextValue $ws.Range("E43") "  +0.97%  "

# Row 44: Maker
Set-TextValue $ws.Range("D44") "1.005.40"
Set-TextValue $ws.Range("E44") "  -7.02%  "

# Row 45: RocketPoolETH
Set-TextValue $ws.Range("D45") "1.806.90"
Set-TextValue $ws.Range("E45") "  +0.03%  "

# Row 46: Aave
Set-TextValue $ws.Range("D46") "56.75"
Set-TextValue $ws.Range("E46") "  -1.58%  "

# Row 47: BabyDogeCoin
Set-TextValue $ws.Range("D47") "0.0₈107"
Set-TextValue $ws.Range("E47") "  -2.53%  "

# Row 48: Frax
Set-TextValue $ws.Range("D48") "0.9987"
Set-TextValue $ws.Range("E48") "  -0.12%  "

# Row 49: EnergySwap
Set-TextValue $ws.Range("D49") "7.973"
Set-TextValue $ws.Range("E49") "  -1.73%  "

# Row 50: Mantle
Set-TextValue $ws.Range("D50") "0.4314"
Set-TextValue $ws.Range("E50") "  +1.90%  "

# Row 51: Cronos
Set-TextValue $ws.Range("E51") "  -0.86%  "
